# Apply "Many updates due to SIMBAD changes" edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Single-cell precision tweaks in column Q (HZ_Detection_Limit) ---
$ws.Range("Q9").Value = 1.616889990960015
$ws.Range("Q29").Value = 2.608985367019665
$ws.Range("Q79").Value = 5.029247845400172
$ws.Range("Q80").Value = 5.06680867658183
$ws.Range("Q92").Value = 33.7642303694855
$ws.Range("Q93").Value = 187.1954334784243

# --- Rows 32-34 get reshuffled (cyclic shift) along with one more
#     precision tweak to the HZ_Detection_Limit of the row that ends
#     up at row 34. Capture the original row values first, then
#     rewrite every cell so the dependent formatting/order matches the
#     target workbook. ---

$origRow32 = @(2443.01, 318753380, 40.179861, 1.199676, 2459148.098617, 15.6692322, 4.562853, 1393.3814802, 2.6868893, 13.2931176, 486.9971964, 8.296900000000001, 23.9258, 4214.44, 4.52845, 0.732115, 2.839885605748058, 0.1818675256571092)
$origRow33 = @(6965.01, 80224448, 103.678273, 24.245141, 2459505.193695, 5.9693397, 3.441038, 105.5424205, 1.3473608, 334.6170929, 1090.8292631, 6.3103, 31.1664, 6007, 4.37729, 1.12993, 2.880011854113967, 0.05173693485062012)
$origRow34 = @(6965.02, 80224448, 103.678273, 24.245141, 2459508.008646, 28.0693949, 6.0393818, 154.3008171, 1.6221344, 42.4755605, 651.1102456, 6.3103, 31.1664, 6007, 4.37729, 1.12993, 2.880011854113967, 0.05173693485062012)

# New row 32 <- old row 33 data (unchanged values)
$newRow32 = $origRow33
# New row 33 <- old row 34 data (unchanged values)
$newRow33 = $origRow34
# New row 34 <- old row 32 data, but HZ_Detection_Limit (col Q, index 16) updates
$newRow34 = $origRow32
$newRow34[16] = 2.970022591198855

for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(32, $col).Value = $newRow32[$col - 1]
    $ws.Cells.Item(33, $col).Value = $newRow33[$col - 1]
    $ws.Cells.Item(34, $col).Value = $newRow34[$col - 1]
}
